# Updated cryptos list on Fri Jun  7 11:58:19 UTC 2024 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: row, new Price (D) value (or $null to leave unchanged), new Volume(1h) (E) value (or $null)
# Price values are prefixed with a leading apostrophe so Excel stores them as text
# (matching the source data) instead of auto-coercing to numbers.
$updates = @(
    @{Row=2;  D="'71.549.16";  E="  +0.73%  "},
    @{Row=3;  D="'3.831.66";   E="  -0.42%  "},
    @{Row=4;  D=$null;         E="  -0.06%  "},
    @{Row=5;  D="'703.74";     E="  -0.17%  "},
    @{Row=6;  D="'172.02";     E="  -0.54%  "},
    @{Row=7;  D="'3.834.93";   E="  -0.35%  "},
    @{Row=8;  D=$null;         E="  -0.05%  "},
    @{Row=9;  D=$null;         E="  +0.19%  "},
    @{Row=10; D=$null;         E="  -1.41%  "},
    @{Row=11; D=$null;         E="  +2.31%  "},
    @{Row=12; D="'0.484";      E="  +5.48%  "},
    @{Row=13; D=$null;         E="  -1.46%  "},
    @{Row=14; D="'36.29";      E="  -1.39%  "},
    @{Row=15; D=$null;         E="  -0.67%  "},
    @{Row=16; D="'3.830.51";   E="  -1.65%  "},
    @{Row=17; D="'71.580.44";  E="  +0.68%  "},
    @{Row=18; D="'7.26";       E="  +0.61%  "},
    @{Row=19; D="'17.57";      E="  +1.17%  "},
    @{Row=20; D=$null;         E="  -0.15%  "},
    @{Row=21; D="'516.58";     E="  +3.89%  "},
    @{Row=22; D="'10.53";      E="  -1.22%  "},
    @{Row=23; D="'0.719";      E="  +0.30%  "},
    @{Row=24; D="'84.17";      E="  -1.39%  "},
    @{Row=25; D=$null;         E="  -2.99%  "},
    @{Row=26; D="'12.80";      E="  +4.84%  "},
    @{Row=28; D="'10.39";      E=$null},
    @{Row=29; D=$null;         E="  +0.06%  "},
    @{Row=31; D=$null;         E="  -5.24%  "},
    @{Row=32; D="'7.41";       E="  -1.65%  "},
    @{Row=34; D="'29.32";      E="  -0.55%  "},
    @{Row=35; D=$null;         E="  -3.51%  "},
    @{Row=36; D=$null;         E="  +0.65%  "},
    @{Row=37; D="'3.786.85";   E="  -0.47%  "},
    @{Row=38; D=$null;         E="  -0.05%  "},
    @{Row=39; D=$null;         E="  -1.84%  "},
    @{Row=40; D="'2.48";       E="  +4.25%  "},
    @{Row=41; D="'6.38";       E="  +5.73%  "},
    @{Row=42; D=$null;         E="  -1.10%  "},
    @{Row=43; D="'3.29";       E="  -2.05%  "},
    @{Row=45; D="'173.34";     E="  +6.00%  "},
    @{Row=46; D=$null;         E="  -0.05%  "},
    @{Row=49; D="'429.05";     E="  +2.28%  "},
    @{Row=50; D=$null;         E="  -0.83%  "},
    @{Row=51; D="'8.60";       E="  -0.12%  "}
)

foreach ($u in $updates) {
    $r = $u.Row
    if ($null -ne $u.D) {
        $ws.Cells.Item($r, 4).Value = $u.D
    }
    if ($null -ne $u.E) {
        $ws.Cells.Item($r, 5).Value = $u.E
    }
}

# Row 47/48 swap: OKB (row 47) and FLOKI (row 48) traded ranking positions,
# each refreshed with a new price/volume reading.
$ws.Cells.Item(47, 2).Value = "FLOKI"
$ws.Cells.Item(47, 3).Value = "https://coinranking.com/coin/fmHk13Rqw+floki-floki"
$ws.Cells.Item(47, 4).Value = "'0.000309"
$ws.Cells.Item(47, 5).Value = "  -5.14%  "

$ws.Cells.Item(48, 2).Value = "OKB"
$ws.Cells.Item(48, 3).Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Cells.Item(48, 4).Value = "'49.95"
$ws.Cells.Item(48, 5).Value = "  +2.78%  "
